$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; remaining columns (B:F) shift left to (A:E)
$ws.Range("A:A").Delete()

# Rename the header "MODEL_CONDITION" -> "MODELCONDITION"
# (this is now in D1 after the column shift)
$ws.Range("D1").Value = "MODELCONDITION"
